# Expert Add Quiz Page Test Cases Added with Assertions
#
# Adds an "assertion" column (E: yes/no) to every existing quiz test-case
# row, adds a new "Recruitment Name" column (F) populated for the first
# data row, and appends two brand-new test-case rows (12 and 13) that
# include a recruitment name / empty-field scenario.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Fix the pre-existing typo/value: row 8's time-limit cell held the
#    literal text "0:00" - it should read "0".
# ---------------------------------------------------------------------
$ws.Range("C8").Value = """0"""

# ---------------------------------------------------------------------
# 2. New "assertion" values (yes/no) for every existing data row.
# ---------------------------------------------------------------------
$ws.Range("E5").Value = "no"
$ws.Range("E6").Value = "no"
$ws.Range("E7").Value = "no"
$ws.Range("E8").Value = "no"
$ws.Range("E2").Value = "yes"
$ws.Range("E3").Value = "yes"
$ws.Range("E4").Value = "yes"
$ws.Range("E9").Value = "yes"
$ws.Range("E10").Value = "yes"
$ws.Range("E11").Value = "yes"

# ---------------------------------------------------------------------
# 3. New header cell for column F (recruitment name) and its first value.
# ---------------------------------------------------------------------
$ws.Range("F1").Value = "Recruitment Name"
$ws.Range("F1").Font.Bold = $true
$ws.Range("F2").Value = "RECRUITMENT18"

# ---------------------------------------------------------------------
# 4. Two brand-new quiz-link test case rows (12 & 13), matching the style
#    of the other "www.google.com" hyperlink rows. Clone row 11's number
#    formatting first (date / time cell formats), then fill the values.
# ---------------------------------------------------------------------
$ws.Range("A11:E11").Copy() | Out-Null
$ws.Range("A12:E12").PasteSpecial(-4122) | Out-Null
$ws.Range("A11:E11").Copy() | Out-Null
$ws.Range("A13:E13").PasteSpecial(-4122) | Out-Null

$ws.Range("A12").Value = "www.google.com"
$ws.Hyperlinks.Add($ws.Range("A12"), "http://www.google.com/")
$ws.Range("A12").Style = "Hyperlink"
$ws.Range("B12").Value = """4/22/2022"""
$ws.Range("C12").Value = """01:00"""
$ws.Range("D12").Value = """5"""
$ws.Range("E12").Value = "no"

$ws.Range("A13").Value = "www.google.com"
$ws.Hyperlinks.Add($ws.Range("A13"), "http://www.google.com/")
$ws.Range("A13").Style = "Hyperlink"
$ws.Range("B13").Value = """4/22/2022"""
$ws.Range("C13").Value = """01:00"""
$ws.Range("D13").Value = """"""
$ws.Range("E13").Value = "yes"

# ---------------------------------------------------------------------
# 5. Column F width, to fit the new "Recruitment Name" header.
# ---------------------------------------------------------------------
$ws.Columns("F").ColumnWidth = 15.166666666666666

# ---------------------------------------------------------------------
# 6. Leave the same selection Excel would land on after typing into F13.
# ---------------------------------------------------------------------
$ws.Range("F13").Select() | Out-Null
